# Agregar prioridad y ocupacion
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J (shifts existing J..O to K..P), then fill header/value
$ws.Columns("J:J").Insert()
$ws.Range("J1").Value = "Prioridad"
$ws.Range("J2").Value = 55

# Append a new "Ocupacion" column at the end (now column Q), copying the
# header style from an existing header cell so it matches (bold/border/center)
$ws.Range("A1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)
$ws.Range("Q1").Value = "Ocupación"
$ws.Range("Q2").Value = 0.9

# Update the rest of row 2 values per the diff
$ws.Range("A2").Value = "vAp0n200s12d1"
$ws.Range("B2").Value = 8312.780000000001
$ws.Range("C2").Value = "Opt"
$ws.Range("D2").Value = 59000
$ws.Range("E2").Value = 59000
$ws.Range("F2").Value = 0

# G2 is a text cell (holds a numeric-looking string) - force text storage
# without permanently changing its format, matching the original inlineStr.
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "25.26"
$ws.Range("G2").Style = "Normal"

$ws.Range("H2").Value = 12.98
$ws.Range("I2").Value = 55

# Columns shifted right by the insert: old K..O values now at L..P, update them
$ws.Range("K2").Value = 53.2
$ws.Range("L2").Value = 37.7
$ws.Range("M2").Value = 5.5
$ws.Range("N2").Value = 3.5
$ws.Range("O2").Value = 9.300000000000001
$ws.Range("P2").Value = 2.8
